$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 to become the "Navel" rate previously on row 8
$ws.Range("A3").Value = "Navel"
$ws.Range("B3").Value = 93
$ws.Range("C3").Value = "All"

# Remove rows 4 through 13 (old Grapefruit .. Valencia rows)
$ws.Range("A4:C13").EntireRow.Delete()
